$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns touched below so that
# numeric-looking strings (e.g. "256.94", "-0.59%") are stored as literal
# text, matching the source workbook's inlineStr cells instead of being
# auto-converted to numbers/percentages by Excel's input parser.
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","E18","E19","D22","E22","D23","E23","D24","E24","D25","E25","E26","D27","E27","D40","E40","B41","C41","D41","E41","B42","C42","D42","E42","B43","C43","D43","E43","B44","C44","D44","E44","D45","E45","E46","D47","E47","E48","E49","E50")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "256.94"
$ws.Range("E2").Value = "-0.59%"
$ws.Range("D3").Value = "27.08"
$ws.Range("E3").Value = "1.02%"
$ws.Range("D4").Value = "4.506"
$ws.Range("E4").Value = "-6.64%"
$ws.Range("D5").Value = "0.05892"
$ws.Range("E5").Value = "-1.33%"
$ws.Range("E6").Value = "-0.86%"
$ws.Range("D7").Value = "0.8532"
$ws.Range("E7").Value = "-2.37%"
$ws.Range("D8").Value = "0.9367"
$ws.Range("E8").Value = "-1.73%"
$ws.Range("D9").Value = "0.1389"
$ws.Range("E9").Value = "-1.88%"
$ws.Range("D10").Value = "0.04760"
$ws.Range("E10").Value = "32.23%"
$ws.Range("D11").Value = "0.07079"
$ws.Range("E11").Value = "-1.87%"
$ws.Range("D12").Value = "0.03074"
$ws.Range("E12").Value = "-2.82%"
$ws.Range("D13").Value = "0.09113"
$ws.Range("E13").Value = "-1.35%"
$ws.Range("D14").Value = "0.001522"
$ws.Range("E14").Value = "-1.89%"
$ws.Range("D15").Value = "0.0006041"
$ws.Range("E15").Value = "-94.33%"
$ws.Range("D16").Value = "0.006044"
$ws.Range("E16").Value = "0.26%"
$ws.Range("D17").Value = "3.493"
$ws.Range("E17").Value = "0.19%"
$ws.Range("E18").Value = "-1.44%"
$ws.Range("E19").Value = "-1.57%"
$ws.Range("D22").Value = "3.917"
$ws.Range("E22").Value = "10.89%"
$ws.Range("D23").Value = "0.04273"
$ws.Range("E23").Value = "1.46%"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").Value = "-0.13%"
$ws.Range("D25").Value = "0.004287"
$ws.Range("E25").Value = "-5.00%"
$ws.Range("E26").Value = "0.01%"
$ws.Range("D27").Value = "0.0001524"
$ws.Range("E27").Value = "2.08%"
$ws.Range("D40").Value = "0.03823"
$ws.Range("E40").Value = "-0.58%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006289"
$ws.Range("E41").Value = "5.20%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1100"
$ws.Range("E42").Value = "-0.30%"
$ws.Range("B43").Value = "LocalTraders"
$ws.Range("C43").Value = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"
$ws.Range("D43").Value = "0.01409"
$ws.Range("E43").Value = "31.81%"
$ws.Range("B44").Value = "CEJI"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D44").Value = "0.001900"
$ws.Range("E44").Value = "-13.63%"
$ws.Range("D45").Value = "0.00005354"
$ws.Range("E45").Value = "-2.50%"
$ws.Range("E46").Value = "0.01%"
$ws.Range("D47").Value = "0.06589"
$ws.Range("E47").Value = "-39.59%"
$ws.Range("E48").Value = "11,738.84%"
$ws.Range("E49").Value = "0.01%"
$ws.Range("E50").Value = "0.01%"
